$wb = $excel.ActiveWorkbook
$wsInv = $wb.Worksheets.Item(1)

# Update selection on Inventory List sheet
$wsInv.Range("F14").Select()

# ---- Gender-List ----
$wsGender = $wb.Worksheets.Add($null, $wsInv)
$wsGender.Name = "Gender-List"
$wsGender.Range("A1").Value = "ID"
$wsGender.Range("B1").Value = "Keterangan"
$wsGender.Range("A2").Value = 1
$wsGender.Range("B2").Value = "Male"
$wsGender.Range("A3").Value = 2
$wsGender.Range("B3").Value = "Female"

# ---- Religion-List ----
$wsReligion = $wb.Worksheets.Add($null, $wsGender)
$wsReligion.Name = "Religion-List"
$wsReligion.Range("A1").Value = "ID"
$wsReligion.Range("B1").Value = "Keterangan"
$wsReligion.Range("A2").Value = 2
$wsReligion.Range("B2").Value = "ISLAM"
$wsReligion.Range("A3").Value = 3
$wsReligion.Range("B3").Value = "Budha"
$wsReligion.Range("A4").Value = 4
$wsReligion.Range("B4").Value = "Katolik"
$wsReligion.Range("A5").Value = 5
$wsReligion.Range("B5").Value = "Hindu"
$wsReligion.Range("A6").Value = 6
$wsReligion.Range("B6").Value = "Kristen"
$wsReligion.Range("A7").Value = 7
$wsReligion.Range("B7").Value = "Kong Hu Chu"
$wsReligion.Range("A8").Value = 8
$wsReligion.Range("B8").Value = "Lain - Lain"

# ---- Job-List ----
$wsJob = $wb.Worksheets.Add($null, $wsReligion)
$wsJob.Name = "Job-List"
$wsJob.Range("A1").Value = "ID"
$wsJob.Range("B1").Value = "Keterangan"
$wsJob.Range("A2").Value = 2
$wsJob.Range("B2").Value = "Pelajar/Mahasiswa"
$wsJob.Range("A3").Value = 3
$wsJob.Range("B3").Value = "Wirausaha"
$wsJob.Range("A4").Value = 4
$wsJob.Range("B4").Value = "TNI"
$wsJob.Range("A5").Value = 5
$wsJob.Range("B5").Value = "KARYAWAN SWASTA"
$wsJob.Range("A6").Value = 6
$wsJob.Range("B6").Value = "Ibu Rumah Tangga"
$wsJob.Range("A7").Value = 8
$wsJob.Range("B7").Value = "Guru"
$wsJob.Range("A8").Value = 9
$wsJob.Range("B8").Value = "Belum Bekerja"
$wsJob.Range("A9").Value = 10
$wsJob.Range("B9").Value = "Pegawai Negri Sipil"
$wsJob.Range("A10").Value = 11
$wsJob.Range("B10").Value = "Pengacara"
$wsJob.Range("A11").Value = 12
$wsJob.Range("B11").Value = "Karyawan BUMN"

# ---- Pendidikan-List ----
$wsPendidikan = $wb.Worksheets.Add($null, $wsJob)
$wsPendidikan.Name = "Pendidikan-List"
$wsPendidikan.Range("A1").Value = "ID"
$wsPendidikan.Range("B1").Value = "Keterangan"
$wsPendidikan.Range("A2").Value = 1
$wsPendidikan.Range("B2").Value = "SD"
$wsPendidikan.Range("A3").Value = 2
$wsPendidikan.Range("B3").Value = "SMPN"
$wsPendidikan.Range("A4").Value = 3
$wsPendidikan.Range("B4").Value = "SMA/SMK"
$wsPendidikan.Range("A5").Value = 4
$wsPendidikan.Range("B5").Value = "D3"
$wsPendidikan.Range("A6").Value = 5
$wsPendidikan.Range("B6").Value = "S1"
$wsPendidikan.Range("A7").Value = 6
$wsPendidikan.Range("B7").Value = "SMP"
$wsPendidikan.Range("A8").Value = 7
$wsPendidikan.Range("B8").Value = "SD"
$wsPendidikan.Range("A9").Value = 8
$wsPendidikan.Range("B9").Value = "TK"
$wsPendidikan.Range("A10").Value = 9
$wsPendidikan.Range("B10").Value = "PAUD"
$wsPendidikan.Range("A11").Value = 10
$wsPendidikan.Range("B11").Value = "D4"
$wsPendidikan.Range("A12").Value = 11
$wsPendidikan.Range("B12").Value = "S2"
$wsPendidikan.Range("A13").Value = 12
$wsPendidikan.Range("B13").Value = "S3"

# Style header rows (gray fill + border) and data rows (border only) for new sheets
foreach ($sheet in @($wsGender, $wsReligion, $wsJob, $wsPendidikan)) {
    $lastRow = $sheet.Range("A1").End(4).Row  # xlDown = 4, but using count below instead
}

$wsGender.Range("A1:B1").Interior.ThemeColor = 1
$wsGender.Range("A1:B1").Interior.TintAndShade = -0.249977111117893
$wsGender.Range("A1:B3").Borders.LineStyle = 1

$wsReligion.Range("A1:B1").Interior.ThemeColor = 1
$wsReligion.Range("A1:B1").Interior.TintAndShade = -0.249977111117893
$wsReligion.Range("A1:B8").Borders.LineStyle = 1

$wsJob.Range("A1:B1").Interior.ThemeColor = 1
$wsJob.Range("A1:B1").Interior.TintAndShade = -0.249977111117893
$wsJob.Range("A1:B11").Borders.LineStyle = 1

$wsPendidikan.Range("A1:B1").Interior.ThemeColor = 1
$wsPendidikan.Range("A1:B1").Interior.TintAndShade = -0.249977111117893
$wsPendidikan.Range("A1:B13").Borders.LineStyle = 1

$wsInv.Select()
